# Append a new row (row 33) to each of the 4 worksheets, copied from row 32
# but with the timestamp in column A advanced by one hour.

$wb = $excel.ActiveWorkbook

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    $srcRow = 32
    $dstRow = 33

    # Read source row values (column A..I => 1..9)
    $vals = @{}
    for ($col = 1; $col -le 9; $col++) {
        $vals[$col] = $ws.Cells.Item($srcRow, $col).Value2
    }

    # Bump the timestamp in column A by one hour (plain string manipulation,
    # values look like "2025-03-05 15:42:06").
    $timeStr = [string]$vals[1]
    $datePart = $timeStr.Substring(0, 10)
    $timePart = $timeStr.Substring(11, 8)
    $pieces = $timePart.Split(":")
    $hour = ([int]$pieces[0] + 1) % 24
    $hourStr = $hour.ToString("00")
    $newTimeStr = $datePart + " " + $hourStr + ":" + $pieces[1] + ":" + $pieces[2]

    # Columns that are stored as numbers in the sheet.
    $numericCols = @(6, 8, 9)

    for ($col = 1; $col -le 9; $col++) {
        $cell = $ws.Cells.Item($dstRow, $col)
        if ($numericCols -contains $col) {
            $cell.NumberFormat = "General"
            $cell.Value = [double]$vals[$col]
        } else {
            # Force text so numeric-looking strings (e.g. the big ID in
            # column G) and the timestamp stay as text, not auto-converted.
            $cell.NumberFormat = "@"
            if ($col -eq 1) {
                $cell.Value = $newTimeStr
            } else {
                $cell.Value = [string]$vals[$col]
            }
        }
    }
}
